##############################################################################
# edit.ps1 - Apply binance_ETHUSDT_data.xlsx update
#
# Commit: "Update with Position Function"
#  * refresh row 157 (last existing row) with corrected OHLCV + indicator
#    values
#  * append 8 new daily rows (158-165, covering 2020-11-04 .. 2020-11-11)
#  * grow the sheet dimension from A1:X157 to A1:X165
##############################################################################

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (matching the workbook convention
# of storing OHLCV numbers as fixed-format strings, e.g. "390.00000000"),
# without leaving a lingering custom number format / cell style behind.
function Set-TextCell($addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

##############################################################################
# Row 157: correct the close/volume/indicator columns (open/low/close_time
# and the GST timestamp label are untouched)
##############################################################################
Set-TextCell "C157" "390.00000000"
Set-TextCell "E157" "387.88000000"
Set-TextCell "F157" "752954.37114000"
Set-TextCell "H157" "286291773.96848450"
$ws.Range("I157").Value = 280412
Set-TextCell "J157" "366506.41600000"
Set-TextCell "K157" "139475874.91715260"
$ws.Range("M157").Value = 387.8799999999997
$ws.Range("N157").Value = 385.4499999999998
$ws.Range("O157").Value = 387.4114285714285
$ws.Range("P157").Value = 393.8559999999999
$ws.Range("Q157").Value = 380.6516666666666
$ws.Range("R157").Value = 387.88
$ws.Range("S157").Value = 387.3761605722846
$ws.Range("T157").Value = 389.3449945121245
$ws.Range("U157").Value = 384.5561140938933
$ws.Range("V157").Value = 4.788880418231201
$ws.Range("W157").Value = 6.733537864916512
$ws.Range("X157").Value = -1.94465744668531

##############################################################################
# Rows 158-165: new daily candles appended after row 157
##############################################################################

# Column A carries the bold/bordered "index" style (same as existing rows);
# copy that formatting onto the new A158:A165 block in one shot so we reuse
# the workbook's existing style instead of minting a new one.
$ws.Range("A157").Copy()
$ws.Range("A158:A165").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---- Row 158 ----
$ws.Range("A158").Value = 156
Set-TextCell "B158" "387.88000000"
Set-TextCell "C158" "408.39000000"
Set-TextCell "D158" "376.47000000"
Set-TextCell "E158" "402.49000000"
Set-TextCell "F158" "1157064.35076000"
Set-TextCell "H158" "454481208.40200980"
Set-TextCell "J158" "600951.59602000"
Set-TextCell "K158" "236190448.28833480"
Set-TextCell "L158" "2020-11-04 08:00:00"
$ws.Range("G158").Value = 1604534399999
$ws.Range("I158").Value = 392310
$ws.Range("M158").Value = 402.4899999999997
$ws.Range("N158").Value = 395.1849999999998
$ws.Range("O158").Value = 389.4485714285714
$ws.Range("P158").Value = 396.1206666666666
$ws.Range("Q158").Value = 382.2806666666665
$ws.Range("R158").Value = 402.49
$ws.Range("S158").Value = 397.4520535240949
$ws.Range("T158").Value = 391.367303048729
$ws.Range("U158").Value = 385.8845575999398
$ws.Range("V158").Value = 5.48274544878916
$ws.Range("W158").Value = 6.483379381691043
$ws.Range("X158").Value = -1.000633932901883

# ---- Row 159 ----
$ws.Range("A159").Value = 157
Set-TextCell "B159" "402.50000000"
Set-TextCell "C159" "420.40000000"
Set-TextCell "D159" "396.14000000"
Set-TextCell "E159" "416.69000000"
Set-TextCell "F159" "1475139.77342000"
Set-TextCell "H159" "601591250.40124900"
Set-TextCell "J159" "733696.12776000"
Set-TextCell "K159" "299466668.01144860"
Set-TextCell "L159" "2020-11-05 08:00:00"
$ws.Range("G159").Value = 1604620799999
$ws.Range("I159").Value = 513574
$ws.Range("M159").Value = 416.6899999999996
$ws.Range("N159").Value = 409.5899999999999
$ws.Range("O159").Value = 393.6714285714285
$ws.Range("P159").Value = 397.8746666666666
$ws.Range("Q159").Value = 384.8119999999999
$ws.Range("R159").Value = 416.69
$ws.Range("S159").Value = 410.2773511746983
$ws.Range("T159").Value = 395.2631025797072
$ws.Range("U159").Value = 388.1664541714699
$ws.Range("V159").Value = 7.096648408237286
$ws.Range("W159").Value = 6.606033187000292
$ws.Range("X159").Value = 0.4906152212369941

# ---- Row 160 ----
$ws.Range("A160").Value = 158
Set-TextCell "B160" "416.73000000"
Set-TextCell "C160" "458.27000000"
Set-TextCell "D160" "414.76000000"
Set-TextCell "E160" "455.91000000"
Set-TextCell "F160" "1682440.96026000"
Set-TextCell "H160" "736015135.62721030"
Set-TextCell "J160" "883105.68936000"
Set-TextCell "K160" "386541087.43089730"
Set-TextCell "L160" "2020-11-06 08:00:00"
$ws.Range("G160").Value = 1604707199999
$ws.Range("I160").Value = 666505
$ws.Range("M160").Value = 455.9099999999997
$ws.Range("N160").Value = 436.2999999999998
$ws.Range("O160").Value = 404.16
$ws.Range("P160").Value = 400.67
$ws.Range("Q160").Value = 388.6229999999999
$ws.Range("R160").Value = 455.91
$ws.Range("S160").Value = 440.6991170582328
$ws.Range("T160").Value = 404.5933944905486
$ws.Range("U160").Value = 393.1845189341669
$ws.Range("V160").Value = 11.4088755563817
$ws.Range("W160").Value = 7.566601660876574
$ws.Range("X160").Value = 3.842273895505127

# ---- Row 161 ----
$ws.Range("A161").Value = 159
Set-TextCell "B161" "455.91000000"
Set-TextCell "C161" "468.28000000"
Set-TextCell "D161" "424.14000000"
Set-TextCell "E161" "435.21000000"
Set-TextCell "F161" "1653416.59610000"
Set-TextCell "H161" "743871175.78212340"
Set-TextCell "J161" "831369.37192000"
Set-TextCell "K161" "374507765.44837650"
Set-TextCell "L161" "2020-11-07 08:00:00"
$ws.Range("G161").Value = 1604793599999
$ws.Range("I161").Value = 674858
$ws.Range("M161").Value = 435.2099999999996
$ws.Range("N161").Value = 445.5599999999998
$ws.Range("O161").Value = 411.1242857142857
$ws.Range("P161").Value = 402.41
$ws.Range("Q161").Value = 391.4316666666665
$ws.Range("R161").Value = 435.21
$ws.Range("S161").Value = 437.0397056860776
$ws.Range("T161").Value = 409.3036414920143
$ws.Range("U161").Value = 396.2975315075429
$ws.Range("V161").Value = 13.00610998447144
$ws.Range("W161").Value = 8.654503325595547
$ws.Range("X161").Value = 4.351606658875889

# ---- Row 162 ----
$ws.Range("A162").Value = 160
Set-TextCell "B162" "435.20000000"
Set-TextCell "C162" "460.10000000"
Set-TextCell "D162" "431.99000000"
Set-TextCell "E162" "454.30000000"
Set-TextCell "F162" "840908.46841000"
Set-TextCell "H162" "375679148.23503960"
Set-TextCell "J162" "424300.36519000"
Set-TextCell "K162" "189644075.40081470"
Set-TextCell "L162" "2020-11-08 08:00:00"
$ws.Range("G162").Value = 1604879999999
$ws.Range("I162").Value = 354341
$ws.Range("M162").Value = 454.2999999999997
$ws.Range("N162").Value = 444.7549999999998
$ws.Range("O162").Value = 419.3571428571428
$ws.Range("P162").Value = 405.2206666666667
$ws.Range("Q162").Value = 394.4076666666664
$ws.Range("R162").Value = 454.3
$ws.Range("S162").Value = 448.5465685620259
$ws.Range("T162").Value = 416.2261581855651
$ws.Range("U162").Value = 400.5940285154788
$ws.Range("V162").Value = 15.63212967008627
$ws.Range("W162").Value = 10.05002859449369
$ws.Range("X162").Value = 5.582101075592579

# ---- Row 163 ----
$ws.Range("A163").Value = 161
Set-TextCell "B163" "454.29000000"
Set-TextCell "C163" "459.10000000"
Set-TextCell "D163" "433.09000000"
Set-TextCell "E163" "444.32000000"
Set-TextCell "F163" "1099213.24641000"
Set-TextCell "H163" "491928164.98124340"
Set-TextCell "J163" "523180.45111000"
Set-TextCell "K163" "234223108.98942740"
Set-TextCell "L163" "2020-11-09 08:00:00"
$ws.Range("G163").Value = 1604966399999
$ws.Range("I163").Value = 469947
$ws.Range("M163").Value = 444.3199999999997
$ws.Range("N163").Value = 449.3099999999997
$ws.Range("O163").Value = 428.1142857142858
$ws.Range("P163").Value = 407.7846666666666
$ws.Range("Q163").Value = 396.8676666666664
$ws.Range("R163").Value = 444.32
$ws.Range("S163").Value = 445.728856187342
$ws.Range("T163").Value = 420.5482876954858
$ws.Range("U163").Value = 403.8330018330985
$ws.Range("V163").Value = 16.71528586238725
$ws.Range("W163").Value = 11.3830800480724
$ws.Range("X163").Value = 5.332205814314843

# ---- Row 164 ----
$ws.Range("A164").Value = 162
Set-TextCell "B164" "444.32000000"
Set-TextCell "C164" "455.00000000"
Set-TextCell "D164" "438.70000000"
Set-TextCell "E164" "450.34000000"
Set-TextCell "F164" "866872.41108000"
Set-TextCell "H164" "388190148.20789510"
Set-TextCell "J164" "437244.02871000"
Set-TextCell "K164" "195856348.18969540"
Set-TextCell "L164" "2020-11-10 08:00:00"
$ws.Range("G164").Value = 1605052799999
$ws.Range("I164").Value = 384228
$ws.Range("M164").Value = 450.3399999999996
$ws.Range("N164").Value = 447.3299999999997
$ws.Range("O164").Value = 437.0371428571429
$ws.Range("P164").Value = 411.6393333333334
$ws.Range("Q164").Value = 399.4069999999998
$ws.Range("R164").Value = 450.34
$ws.Range("S164").Value = 448.8029520624473
$ws.Range("T164").Value = 425.1316280500333
$ws.Range("U164").Value = 407.2779769378895
$ws.Range("V164").Value = 17.85365111214378
$ws.Range("W164").Value = 12.67719426088668
$ws.Range("X164").Value = 5.176456851257099

# ---- Row 165 ----
$ws.Range("A165").Value = 163
Set-TextCell "B165" "450.34000000"
Set-TextCell "C165" "474.00000000"
Set-TextCell "D165" "449.28000000"
Set-TextCell "E165" "469.23000000"
Set-TextCell "F165" "872086.75803000"
Set-TextCell "H165" "402974529.15965530"
Set-TextCell "J165" "461768.48127000"
Set-TextCell "K165" "213396558.17390210"
Set-TextCell "L165" "2020-11-11 08:00:00"
$ws.Range("G165").Value = 1605139199999
$ws.Range("I165").Value = 392598
$ws.Range("M165").Value = 469.2299999999997
$ws.Range("N165").Value = 459.7849999999997
$ws.Range("O165").Value = 446.5714285714286
$ws.Range("P165").Value = 416.0246666666666
$ws.Range("Q165").Value = 402.1643333333331
$ws.Range("R165").Value = 469.23
$ws.Range("S165").Value = 462.4209840208158
$ws.Range("T165").Value = 431.9159929654214
$ws.Range("U165").Value = 411.8670308267658
$ws.Range("V165").Value = 20.04896213865561
$ws.Range("W165").Value = 14.15154783644047
$ws.Range("X165").Value = 5.89741430221514
